$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 254
$ws.Cells.Item(3, 6).Value = 2655
$ws.Cells.Item(6, 6).Value = 36
$ws.Cells.Item(7, 6).Value = 2055
$ws.Cells.Item(8, 6).Value = 1782
$ws.Cells.Item(11, 6).Value = 2442
$ws.Cells.Item(12, 6).Value = 533
$ws.Cells.Item(13, 6).Value = 217
$ws.Cells.Item(16, 6).Value = 117
$ws.Cells.Item(18, 6).Value = 9028
$ws.Cells.Item(19, 6).Value = 54
$ws.Cells.Item(20, 6).Value = 7001
$ws.Cells.Item(21, 6).Value = 11441
$ws.Cells.Item(22, 6).Value = 125
$ws.Cells.Item(24, 6).Value = 227
$ws.Cells.Item(25, 6).Value = 321
$ws.Cells.Item(27, 6).Value = 2519
$ws.Cells.Item(30, 6).Value = 2442
$ws.Cells.Item(31, 6).Value = 625
$ws.Cells.Item(32, 6).Value = 41
$ws.Cells.Item(33, 6).Value = 4485
$ws.Cells.Item(34, 6).Value = 800
$ws.Cells.Item(35, 6).Value = 337
$ws.Cells.Item(36, 6).Value = 35
$ws.Cells.Item(37, 6).Value = 497

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 67
$ws.Cells.Item(4, 6).Value = 14
$ws.Cells.Item(14, 6).Value = 63
$ws.Cells.Item(16, 6).Value = 96

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 791
$ws.Cells.Item(3, 6).Value = 623

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 791
$ws.Cells.Item(3, 6).Value = 623
$ws.Cells.Item(4, 6).Value = 254
$ws.Cells.Item(6, 6).Value = 2655
$ws.Cells.Item(7, 6).Value = 67
$ws.Cells.Item(9, 6).Value = 36
$ws.Cells.Item(10, 6).Value = 2055
$ws.Cells.Item(11, 6).Value = 14
$ws.Cells.Item(12, 6).Value = 1782
$ws.Cells.Item(15, 6).Value = 2442
$ws.Cells.Item(17, 6).Value = 533
$ws.Cells.Item(18, 6).Value = 217
$ws.Cells.Item(21, 6).Value = 117
$ws.Cells.Item(23, 6).Value = 9028
$ws.Cells.Item(24, 6).Value = 54
$ws.Cells.Item(25, 6).Value = 7002
$ws.Cells.Item(26, 6).Value = 11441
$ws.Cells.Item(28, 6).Value = 125
$ws.Cells.Item(29, 6).Value = 227
$ws.Cells.Item(30, 6).Value = 321
$ws.Cells.Item(34, 6).Value = 2519
$ws.Cells.Item(39, 6).Value = 41
$ws.Cells.Item(40, 6).Value = 4485
$ws.Cells.Item(41, 6).Value = 63
$ws.Cells.Item(43, 6).Value = 96
$ws.Cells.Item(46, 6).Value = 497
